$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 35.07914438502674
$ws.Range("I2").Value = 36.52299465240642
$ws.Range("J2").Value = 37.76577540106952
$ws.Range("K2").Value = 38.80748663101604
$ws.Range("L2").Value = 39.64812834224599
$ws.Range("M2").Value = 40.28770053475936
$ws.Range("N2").Value = 40.72620320855615
$ws.Range("O2").Value = 40.96363636363637

$ws.Range("H3").Value = 16.51764705882353
$ws.Range("I3").Value = 17.28235294117647
$ws.Range("J3").Value = 17.95294117647059
$ws.Range("K3").Value = 18.52941176470588
$ws.Range("L3").Value = 19.01176470588235
$ws.Range("M3").Value = 19.4
$ws.Range("N3").Value = 19.69411764705882
$ws.Range("O3").Value = 19.89411764705882

$ws.Range("H4").Value = 18.56149732620321
$ws.Range("I4").Value = 19.24064171122994
$ws.Range("J4").Value = 19.81283422459893
$ws.Range("K4").Value = 20.27807486631016
$ws.Range("L4").Value = 20.63636363636364
$ws.Range("M4").Value = 20.88770053475936
$ws.Range("N4").Value = 21.03208556149733
$ws.Range("O4").Value = 21.06951871657754

$ws.Range("H5").Value = 35.07914438502674
$ws.Range("I5").Value = 36.52299465240642
$ws.Range("J5").Value = 37.76577540106952
$ws.Range("K5").Value = 38.80748663101604
$ws.Range("L5").Value = 39.64812834224599
$ws.Range("M5").Value = 40.28770053475936
$ws.Range("N5").Value = 40.72620320855615
$ws.Range("O5").Value = 40.96363636363637

$ws.Range("H6").Value = 16.51764705882353
$ws.Range("I6").Value = 17.28235294117647
$ws.Range("J6").Value = 17.95294117647059
$ws.Range("K6").Value = 18.52941176470588
$ws.Range("L6").Value = 19.01176470588235
$ws.Range("M6").Value = 19.4
$ws.Range("N6").Value = 19.69411764705882
$ws.Range("O6").Value = 19.89411764705882

$ws.Range("H7").Value = 18.56149732620321
$ws.Range("I7").Value = 19.24064171122994
$ws.Range("J7").Value = 19.81283422459893
$ws.Range("K7").Value = 20.27807486631016
$ws.Range("L7").Value = 20.63636363636364
$ws.Range("M7").Value = 20.88770053475936
$ws.Range("N7").Value = 21.03208556149733
$ws.Range("O7").Value = 21.06951871657754
